$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test_cases")

# Remove the "profile page" test case row (row 10), which also removes the
# now-unused shared strings ("profile page", "profile page test",
# "Editing the profile", "test_changing_profile").
$ws.Rows.Item(10).Delete() | Out-Null

# Leave the active selection on the (now empty) row just below the data,
# matching where the deleted row used to be.
$ws.Range("A10:H10").Select() | Out-Null
